# Update the servo calibration "Tables" worksheet with new F-column values
# and move the active selection to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")

# New values for column F, rows 4-22 (min pulse width calibration column)
$fValues = @{
    4  = 680
    5  = 770
    6  = 855
    7  = 930
    8  = 1030
    9  = 1115
    10 = 1210
    11 = 1300
    12 = 1400
    13 = 1500
    14 = 1610
    15 = 1700
    16 = 1790
    17 = 1890
    18 = 1970
    19 = 2075
    20 = 2180
    21 = 2280
    22 = 2380
}

foreach ($row in $fValues.Keys) {
    $ws.Range("F$row").Value = $fValues[$row]
}

# Recalculate the workbook so dependent formulas (hex lookup sheets) update
$excel.CalculateFullRebuild()

# Update the stored selection to match the edited range
$ws.Activate()
$ws.Range("F4:F22").Select()
